# Updated: st 14. 01. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revisions to previously reported daily figures (AgTests / AgPosit columns)
$ws.Range("H289").Value = 65124
$ws.Range("I289").Value = 3753

$ws.Range("H293").Value = 83032
$ws.Range("I293").Value = 5864

$ws.Range("H294").Value = 92164
$ws.Range("I294").Value = 5104

$ws.Range("H299").Value = 65417
$ws.Range("I299").Value = 6852

$ws.Range("H300").Value = 71006
$ws.Range("I300").Value = 6956

$ws.Range("H301").Value = 69975
$ws.Range("I301").Value = 5553

$ws.Range("H302").Value = 73049
$ws.Range("I302").Value = 5321

$ws.Range("H306").Value = 70612
$ws.Range("I306").Value = 7170

$ws.Range("H307").Value = 72611
$ws.Range("I307").Value = 6265

$ws.Range("H309").Value = 56861
$ws.Range("I309").Value = 3950

$ws.Range("H310").Value = 89629
$ws.Range("I310").Value = 5362

$ws.Range("H311").Value = 37383
$ws.Range("I311").Value = 1400

$ws.Range("H312").Value = 40569
$ws.Range("I312").Value = 1225

$ws.Range("H313").Value = 71763
$ws.Range("I313").Value = 3476

$ws.Range("H314").Value = 63716
$ws.Range("I314").Value = 3279

# Append new daily record for 2021-01-13 (serial date 44209)
$ws.Range("A315").Value = 44209
$ws.Range("B315").Value = 217978
$ws.Range("C315").Value = 160401
$ws.Range("D315").Value = 54317
$ws.Range("E315").Value = 12342
$ws.Range("F315").Value = 2923
$ws.Range("G315").Value = 3260
$ws.Range("H315").Value = 63435
$ws.Range("I315").Value = 2978
